# Refresh the cryptocurrency price (D) and 1h volume change (E) columns
# on Sheet1 to the latest scraped values (GitHub Actions cron update).
#
# Column D holds price strings such as "71.935.11" or "1.00". Most of
# them contain two '.' separators (thousands + decimal) so Excel leaves
# them as text automatically, but a handful look like plain decimal
# numbers (e.g. "1.00", "589.04"). Typing those into a General-formatted
# cell would make Excel reinterpret them as numbers ("1.00" -> 1, losing
# the trailing zero), so -- exactly like a human typing in Excel -- we
# prefix those with a leading apostrophe to force literal text entry and
# keep the digits exactly as scraped.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '71.935.11'
$ws.Range("E2").Value = '  +3.09%  '
$ws.Range("D3").Value = '3.645.38'
$ws.Range("E3").Value = '  +6.81%  '
$ws.Range("D4").Value = '''1.00'
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").Value = '''589.04'
$ws.Range("E5").Value = '  +0.22%  '
$ws.Range("D6").Value = '''181.09'
$ws.Range("E6").Value = '  -0.53%  '
$ws.Range("D7").Value = '3.635.07'
$ws.Range("E7").Value = '  +6.70%  '
$ws.Range("D8").Value = '''0.617'
$ws.Range("E8").Value = '  +2.65%  '
$ws.Range("E9").Value = '  +0.04%  '
$ws.Range("D10").Value = '''0.203'
$ws.Range("E10").Value = '  -0.23%  '
$ws.Range("E11").Value = '  +2.41%  '
$ws.Range("D12").Value = '''49.95'
$ws.Range("E12").Value = '  +2.60%  '
$ws.Range("E13").Value = '  -0.72%  '
$ws.Range("D14").Value = '''685.41'
$ws.Range("E14").Value = '  -0.66%  '
$ws.Range("D15").Value = '4.234.00'
$ws.Range("E15").Value = '  +6.79%  '
$ws.Range("D16").Value = '''9.06'
$ws.Range("E16").Value = '  +3.82%  '
$ws.Range("D17").Value = '3.668.29'
$ws.Range("E17").Value = '  +7.69%  '
$ws.Range("D18").Value = '71.950.95'
$ws.Range("E18").Value = '  +3.11%  '
$ws.Range("E19").Value = '  +1.63%  '
$ws.Range("D20").Value = '''18.36'
$ws.Range("E20").Value = '  +3.02%  '
$ws.Range("D21").Value = '''11.66'
$ws.Range("E21").Value = '  +2.08%  '
$ws.Range("D22").Value = '''0.941'
$ws.Range("E22").Value = '  +2.71%  '
$ws.Range("D23").Value = '''5.92'
$ws.Range("E23").Value = '  +9.97%  '
$ws.Range("D24").Value = '''17.88'
$ws.Range("E24").Value = '  +2.82%  '
$ws.Range("D25").Value = '''103.42'
$ws.Range("E25").Value = '  -0.17%  '
$ws.Range("E26").Value = '  +1.60%  '
$ws.Range("E27").Value = '  +4.68%  '
$ws.Range("D28").Value = '''10.03'
$ws.Range("E28").Value = '  +2.50%  '
$ws.Range("D29").Value = '''35.13'
$ws.Range("E29").Value = '  +2.82%  '
$ws.Range("D30").Value = '''9.27'
$ws.Range("E30").Value = '  +4.40%  '
$ws.Range("D31").Value = '''7.36'
$ws.Range("E31").Value = '  +5.15%  '
$ws.Range("D32").Value = '''4.18'
$ws.Range("E32").Value = '  +12.89%  '
$ws.Range("D33").Value = '''592.53'
$ws.Range("E33").Value = '  +6.14%  '
$ws.Range("D34").Value = '''11.36'
$ws.Range("E34").Value = '  +1.48%  '
$ws.Range("E35").Value = '  +1.71%  '
$ws.Range("D36").Value = '''59.47'
$ws.Range("E36").Value = '  +1.31%  '
$ws.Range("E37").Value = '  -0.01%  '
$ws.Range("D38").Value = '3.687.94'
$ws.Range("E38").Value = '  +0.39%  '
$ws.Range("E39").Value = '  +0.78%  '
$ws.Range("D40").Value = '''35.79'
$ws.Range("E40").Value = '  -0.95%  '
$ws.Range("D41").Value = '0.0₃0767'
$ws.Range("E41").Value = '  +2.99%  '
$ws.Range("D42").Value = '''0.0472'
$ws.Range("E42").Value = '  +9.62%  '
$ws.Range("D43").Value = '''3.43'
$ws.Range("E43").Value = '  +4.46%  '
$ws.Range("D44").Value = '''2.78'
$ws.Range("E44").Value = '  +2.76%  '
$ws.Range("E45").Value = '  +2.00%  '
$ws.Range("D46").Value = '''3.41'
$ws.Range("E46").Value = '  +1.79%  '
$ws.Range("D47").Value = '''2.82'
$ws.Range("E47").Value = '  +4.82%  '
$ws.Range("E48").Value = '  +2.81%  '
$ws.Range("E49").Value = '  +3.68%  '
$ws.Range("E50").Value = '  -0.14%  '
$ws.Range("D51").Value = '''132.11'
$ws.Range("E51").Value = '  +0.79%  '
